$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)

# Remove the "Straight Connector 111" (id=112) connector shape that was
# deleted as part of the Phase 1 workflow update.
foreach ($sh in $s.Shapes) {
    if ($sh.Name -eq "Straight Connector 111") {
        $sh.Delete()
        break
    }
}
